$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Simplify the "Lab 3" (11-Introduction To Testing / 12-Accessing databases) block:
# Add the missing "10-Interfaces" Lab entry as a new row 33, and shorten the
# "11-Introduction To Testing" Teach/Lab durations.

# New row 33: 10-Interfaces / 45 / Lab
$ws.Range("A33").Value = "10-Interfaces"
$ws.Range("B33").Value = 45
$ws.Range("C33").Value = "Lab"

# Shorten existing Teach/Lab durations for 11-Introduction To Testing
$ws.Range("B34").Value = 35
$ws.Range("B35").Value = 45

# Update the day total to include the new row
$ws.Range("B39").Formula = "=SUM(B33:B38)"
$ws.Range("C39").Formula = "=B39/60"

# Extend the totals row formatting (numFmtId 2 "0.00") across A:F
$ws.Range("A39").NumberFormat = "0.00"
$ws.Range("B39").NumberFormat = "0.00"
$ws.Range("C39").NumberFormat = "0.00"
$ws.Range("D39").NumberFormat = "0.00"
$ws.Range("E39").NumberFormat = "0.00"
$ws.Range("F39").NumberFormat = "0.00"

# Reset the view: scroll back to the top and select B1 instead of the old C39
$null = $ws.Range("B1").Select()
